$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '56.911.26'
$ws.Cells.Item(2, 5).Value = '  +0.24%  '

$ws.Cells.Item(3, 4).Value = '2.428.68'
$ws.Cells.Item(3, 5).Value = '  -3.05%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.26%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '487.27'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.74%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '153.02'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.65%  '

$ws.Cells.Item(7, 5).Value = '  +0.30%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.604'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +17.36%  '

$ws.Cells.Item(9, 4).Value = '2.425.20'
$ws.Cells.Item(9, 5).Value = '  -3.94%  '

$ws.Cells.Item(10, 5).Value = '  -0.15%  '

$ws.Cells.Item(11, 5).Value = '  -1.20%  '

$ws.Cells.Item(12, 5).Value = '  -0.88%  '

$ws.Cells.Item(14, 4).Value = '2.848.54'
$ws.Cells.Item(14, 5).Value = '  -3.16%  '

$ws.Cells.Item(15, 4).Value = '57.062.50'
$ws.Cells.Item(15, 5).Value = '  +0.21%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '20.75'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -3.99%  '

$ws.Cells.Item(17, 5).Value = '  -3.11%  '

$ws.Cells.Item(18, 4).Value = '2.430.24'
$ws.Cells.Item(18, 5).Value = '  -3.23%  '

$ws.Cells.Item(19, 5).Value = '  +3.80%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '324.07'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -0.44%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '9.97'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -4.10%  '

$ws.Cells.Item(22, 5).Value = '  +0.28%  '

$ws.Cells.Item(23, 5).Value = '  +0.17%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '58.11'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -1.89%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.408'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -1.43%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +0.00%  '

$ws.Cells.Item(27, 5).Value = '  -2.66%  '

$ws.Cells.Item(28, 4).Value = '2.532.49'
$ws.Cells.Item(28, 5).Value = '  -3.10%  '

$ws.Cells.Item(29, 5).Value = '  -4.99%  '

$ws.Cells.Item(30, 4).Value = '0.0₃0782'
$ws.Cells.Item(30, 5).Value = '  -5.20%  '

$ws.Cells.Item(31, 5).Value = '  +0.21%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '149.83'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -1.84%  '

$ws.Cells.Item(33, 5).Value = '  +0.28%  '

$ws.Cells.Item(34, 5).Value = '  -0.75%  '

$ws.Cells.Item(35, 5).Value = '  +1.21%  '

$ws.Cells.Item(36, 5).Value = '  -1.76%  '

$ws.Cells.Item(37, 5).Value = '  -3.04%  '

$ws.Cells.Item(38, 5).Value = '  -3.82%  '

$ws.Cells.Item(39, 5).Value = '  +8.55%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '34.08'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.81%  '

$ws.Cells.Item(41, 5).Value = '  -0.16%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.37'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -2.56%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.32%  '

$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.592'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -5.25%  '

$ws.Cells.Item(45, 2).Value = 'Bittensor'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '268.24'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.14%  '

$ws.Cells.Item(46, 5).Value = '  -6.46%  '

$ws.Cells.Item(47, 5).Value = '  -0.05%  '

$ws.Cells.Item(48, 5).Value = '  -1.70%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '4.58'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -8.34%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '17.39'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -3.34%  '

$ws.Cells.Item(51, 4).Value = '1.868.69'
$ws.Cells.Item(51, 5).Value = '  -2.32%  '
